# Atualizei dados add e bibi
# Refreshes the computed "previsao_retorno" figures for a handful of
# clients whose "dias sem comprar" (and a couple of other derived
# metrics) moved because the report was regenerated on a later date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo_por_Cliente")

# --- situacao (coluna J) : meses sem comprar atualizados ---------------
$ws.Range("J5").Value = "INATIVO - 15.0 meses sem comprar"
$ws.Range("J6").Value = "INATIVO - 16.5 meses sem comprar"
$ws.Range("J16").Value = "INATIVO - 40.0 meses sem comprar"
$ws.Range("J24").Value = "INATIVO - 38.0 meses sem comprar"
$ws.Range("J30").Value = "INATIVO - 7.0 meses sem comprar"
$ws.Range("J38").Value = "INATIVO - 32.4 meses sem comprar"
$ws.Range("J44").Value = "INATIVO - 6.3 meses sem comprar"
$ws.Range("J45").Value = "INATIVO - 16.1 meses sem comprar"
$ws.Range("J49").Value = "INATIVO - 7.7 meses sem comprar"
$ws.Range("J64").Value = "INATIVO - 28.0 meses sem comprar"
$ws.Range("J68").Value = "INATIVO - 11.5 meses sem comprar"
$ws.Range("J75").Value = "INATIVO - 7.9 meses sem comprar"
$ws.Range("J78").Value = "INATIVO - 6.3 meses sem comprar"
$ws.Range("J87").Value = "INATIVO - 15.2 meses sem comprar"
$ws.Range("J89").Value = "INATIVO - 11.9 meses sem comprar"
$ws.Range("J90").Value = "INATIVO - 11.3 meses sem comprar"
$ws.Range("J97").Value = "INATIVO - 33.1 meses sem comprar"
$ws.Range("J101").Value = "INATIVO - 37.4 meses sem comprar"
$ws.Range("J103").Value = "INATIVO - 14.7 meses sem comprar"

# --- linha 84 : cliente ficou inativo (proxima_compra vira rotulo) -----
$ws.Range("I84").NumberFormat = "dd/mm/yyyy"
$ws.Range("I84").Value = "INATIVO"
$ws.Range("J84").Value = "INATIVO - 6.0 meses sem comprar"

# --- linha 113 : nova compra registrada, metricas recalculadas ---------
$ws.Range("B113").Value = 0.92
$ws.Range("C113").Value = 0.83
$ws.Range("E113").Value = 16239
$ws.Range("H113").Value = 45840.66871527778
$ws.Range("I113").Value = 45841.66871527778
